# Actualización 11 de Mayo - Mañana
# Updates the "Rescatables" sheet: refreshes the 3 existing student rows
# with corrected data and appends 4 additional rescatable students.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Matricula (Mat) numbers for rows 2-8, in final row order.
$mat = @{
    2 = 20330051920337
    3 = 20330051920339
    4 = 20330051920383
    5 = 20330051920335
    6 = 20330051920336
    7 = 20330051920322
    8 = 18330051920357
}

# Apellido Paterno
$paterno = @{
    2 = "BARRAGAN"
    3 = "CORTES"
    4 = "XOTLANIHUA"
    5 = "VAZQUEZ"
    6 = "BACILIO"
    7 = "GARCIA"
    8 = "BRAVO"
}

# Apellido Materno
$materno = @{
    2 = "VILLALBA"
    3 = "CONTRERAS"
    4 = "MARTINEZ"
    5 = "TZIZIHUA"
    6 = "ATILANO"
    7 = "FLORES"
    8 = "REYES"
}

# Nombres
$nombres = @{
    2 = "ADAN"
    3 = "GERARDO"
    4 = "ANGELA MONTSERRAT"
    5 = "DORA LUZ"
    6 = "YAIR"
    7 = "MARCOS"
    8 = "PATRICIA MARLENE"
}

# Nombre_Largo (materia rescatable)
$nombreLargo = @{
    2 = "APLICA ESTRUCTURAS DE CONTROL CON UN LENGUAJE DE PROGRAMACIÓN"
    3 = "APLICA ESTRUCTURAS DE CONTROL CON UN LENGUAJE DE PROGRAMACIÓN"
    4 = "APLICA ESTRUCTURAS DE CONTROL CON UN LENGUAJE DE PROGRAMACIÓN"
    5 = "ENSAMBLA E INSTALA CONTROLADORES Y DISPOSITIVOS PERIFÉRICOS"
    6 = "APLICA ESTRUCTURAS DE CONTROL CON UN LENGUAJE DE PROGRAMACIÓN"
    7 = "ENSAMBLA E INSTALA CONTROLADORES Y DISPOSITIVOS PERIFÉRICOS"
    8 = "DESARROLLA APLICACIONES MÓVILES PARA ANDROID"
}

# Grupo (matching materia)
$grupo = @{
    2 = "2APM"
    3 = "2APM"
    4 = "2APM"
    5 = "2ASV"
    6 = "2APM"
    7 = "2ASV"
    8 = "6APM"
}

# Reprobadas count
$reprobadas = @{
    2 = 2
    3 = 2
    4 = 2
    5 = 2
    6 = 1
    7 = 1
    8 = 1
}

# Write column by column (A, then B, then C, then D, then E, then F, then G)
# across all rows so new entries land in the shared-string table in the
# same sequence the workbook ends up with.
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 1).Value = $mat[$r]
}
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 2).Value = $paterno[$r]
}
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 3).Value = $materno[$r]
}
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 4).Value = $nombres[$r]
}
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 5).Value = $nombreLargo[$r]
}
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 6).Value = $grupo[$r]
}
foreach ($r in 2..8) {
    $ws.Cells.Item($r, 7).Value = $reprobadas[$r]
}
